$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 114, shifting existing rows 114..168 down to 115..169.
$ws.Rows.Item(114).Insert()

# Populate the newly inserted row 114 with the new record's data.
$ws.Cells.Item(114, 1).Value2 = 9
$ws.Cells.Item(114, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(114, 3).Value2 = "Metropolitana"
$ws.Range("D114").Value2 = 45097
$ws.Cells.Item(114, 5).Value2 = 13
$ws.Cells.Item(114, 6).Value2 = 100112022
$ws.Cells.Item(114, 7).Value2 = "Arveja Verde"
$ws.Cells.Item(114, 8).Value2 = "Perfection"
$ws.Cells.Item(114, 9).Value2 = "Primera"
$ws.Cells.Item(114, 10).Value2 = 43
$ws.Cells.Item(114, 11).Value2 = 36000
$ws.Cells.Item(114, 12).Value2 = 38000
$ws.Cells.Item(114, 13).Value2 = 37023
$ws.Cells.Item(114, 14).Value2 = "`$/saco 25 kilos"
$ws.Cells.Item(114, 15).Value2 = "Provincia de Limarí"
$ws.Cells.Item(114, 16).Value2 = 1481
$ws.Cells.Item(114, 17).Value2 = 25
$ws.Cells.Item(114, 18).Value2 = "Hortaliza"
